# "update data with resort sheetname"
# The workbook originally has two sheets, in tab order:
#   1) "2021-Q3"  (fund-holding detail table, A1:H30)
#   2) "总计"      (quarter summary table, A1:D2)
# This edit re-sorts the sheet tabs so the summary sheet "总计" comes first,
# followed by "2021-Q3" - i.e. it moves "总计" to be the first tab.

$wb = $excel.ActiveWorkbook

$summarySheet = $wb.Worksheets.Item("总计")
$firstSheet   = $wb.Worksheets.Item(1)

# Move "总计" so that it sits right before whatever sheet is currently first.
$summarySheet.Move($firstSheet)
